$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the right of the existing "sum" column (G),
# reusing the same bold/centered/bordered header style as the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill the new column's data rows with 0 (no special style, like the other
# numeric data cells).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
